# Apply cryptos list update (Thu Jul 20 21:48:14 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.810.96"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.892.45"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7982"
$ws.Range("E5").Value = "  -3.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.90"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9989"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3164"
$ws.Range("E8").Value = "  -2.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.51"
$ws.Range("E9").Value = "  -4.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07047"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7676"
$ws.Range("E12").Value = "  +2.64%  "
$ws.Range("D13").Value = "1.895.17"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.297"
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.22"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "29.820.01"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.934"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.01"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007710"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.188"
$ws.Range("E21").Value = "  +18.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "2.135.70"
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9990"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1669"
$ws.Range("E25").Value = "  +4.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.318"
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.79"
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.69"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.058"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.393"
$ws.Range("E30").Value = "  +1.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.536"
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.415"
$ws.Range("E32").Value = "  +3.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05642"
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.051"
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.263"
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7400"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.000"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.622"
$ws.Range("E38").Value = "  -3.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01908"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.775"
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4417"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.60"
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("E43").Value = "  -2.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8450"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9984"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.029.34"
$ws.Range("E46").Value = "  +3.95%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.56"
$ws.Range("E47").Value = "  +1.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.873"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.926"
$ws.Range("E49").Value = "  +2.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.427"
$ws.Range("E50").Value = "  -2.21%  "
$ws.Range("D51").Value = "2.031.69"
$ws.Range("E51").Value = "  -0.76%  "
